$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume symbol-list update (GitHub Actions scrape refresh).
# Each entry is "CellRef|NewValue" - values are written as literal text so the
# sheet keeps matching its original text-formatted Price/Volume(1h) columns
# (prices like "304.13" and percentages like "0.30%" stored as strings, not
# auto-converted numbers/percentages).
$updates = @(
    "D2|304.13",
    "E2|0.30%",
    "D3|37.20",
    "E3|3.71%",
    "D4|5.035",
    "E4|-2.55%",
    "D5|0.07857",
    "E5|-0.09%",
    "D6|2.209",
    "E6|-4.56%",
    "D7|7.988",
    "E7|-0.86%",
    "D8|0.9262",
    "E8|0.03%",
    "D9|0.09832",
    "E9|-2.99%",
    "D10|0.1879",
    "E10|2.63%",
    "D11|0.08644",
    "E11|1.27%",
    "D12|0.03698",
    "E12|9.13%",
    "D13|0.09913",
    "E13|-0.09%",
    "D14|0.001470",
    "E14|-0.20%",
    "D15|0.005687",
    "E15|-1.10%",
    "D16|3.468",
    "E16|-0.43%",
    "D17|4.012",
    "E17|1.06%",
    "D18|2.252",
    "E18|5.95%",
    "D19|0.3409",
    "E19|-0.66%",
    "D20|0.1300",
    "E20|-1.81%",
    "D21|4.769",
    "E21|5.22%",
    "E22|-0.70%",
    "D23|0.04625",
    "E23|0.01%",
    "E24|3.32%",
    "D25|0.004486",
    "E25|-0.10%",
    "D26|0.0001405",
    "E26|8.56%",
    "D27|0.0002718",
    "E27|-19.65%",
    "D39|0.01833",
    "E39|5.01%",
    "D40|0.04773",
    "E40|0.87%",
    "D41|0.008062",
    "E41|2.89%",
    "D42|0.1406",
    "E42|-0.68%",
    "D43|0.007563",
    "E43|-13.82%",
    "D44|0.002113",
    "E44|-4.33%",
    "D45|0.01043",
    "E45|13.98%",
    "D46|0.00006281",
    "E46|4.14%",
    "D47|0.00000000753",
    "E47|0.79%",
    "D48|0.0005805",
    "E48|0.07%",
    "D49|30.71",
    "E49|429.08%",
    "D50|0.002690",
    "E50|0.46%",
    "D51|0.00002107",
    "E51|0.79%"
)

foreach ($entry in $updates) {
    $parts = $entry.Split("|")
    $cellRef = $parts[0]
    $newValue = $parts[1]

    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}
